$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 570.8889
$ws.Range("J28").Value = 625
$ws.Range("L28").Value = 625
$ws.Range("N28").Value = -1595
$ws.Range("H80").Value = 2250
$ws.Range("I80").Value = 2250
$ws.Range("K80").Value = 6750
$ws.Range("M80").Value = -5752
$ws.Range("H83").Value = 2250
$ws.Range("I83").Value = 2250
$ws.Range("K83").Value = 20250
$ws.Range("M83").Value = -15258
$ws.Range("H137").Value = 2199.625
$ws.Range("I137").Value = 1849.8334
$ws.Range("J137").Value = 3249
$ws.Range("K137").Value = 5549.5002
$ws.Range("L137").Value = 9747
$ws.Range("M137").Value = -2999.5002
$ws.Range("N137").Value = -14847
$ws.Range("H138").Value = 8057.926
$ws.Range("J138").Value = 8214
$ws.Range("L138").Value = 24642
$ws.Range("N138").Value = -34922

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1986.6428
$ws.Range("I45").Value = 1986.6428
$ws.Range("K45").Value = 1986.6428
$ws.Range("M45").Value = -1609.6428
$ws.Range("H88").Value = 2822.5
$ws.Range("I88").Value = 2330
$ws.Range("K88").Value = 2330
$ws.Range("M88").Value = -1924
$ws.Range("H91").Value = 2822.5
$ws.Range("I91").Value = 2330
$ws.Range("K91").Value = 2330
$ws.Range("M91").Value = -926
$ws.Range("H102").Value = 1194.3334
$ws.Range("I102").Value = 1183.4
$ws.Range("J102").Value = 1249
$ws.Range("K102").Value = 1183.4
$ws.Range("L102").Value = 1249
$ws.Range("M102").Value = 438.5999999999999
$ws.Range("N102").Value = -4493
$ws.Range("H122").Value = 16997.6
$ws.Range("I122").Value = 19997
$ws.Range("K122").Value = 59991
$ws.Range("M122").Value = -57541

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 60074
$ws.Range("J35").Value = 60074
$ws.Range("L35").Value = 60074
$ws.Range("N35").Value = -60694
$ws.Range("H86").Value = 2362.3684
$ws.Range("I86").Value = 2537.3125
$ws.Range("K86").Value = 2537.3125
$ws.Range("M86").Value = -1414.3125
$ws.Range("H89").Value = 2362.3684
$ws.Range("I89").Value = 2537.3125
$ws.Range("K89").Value = 12686.5625
$ws.Range("M89").Value = -7070.5625
$ws.Range("H94").Value = 499.66666
$ws.Range("J94").Value = 490
$ws.Range("L94").Value = 490
$ws.Range("N94").Value = -1392
$ws.Range("H99").Value = 397
$ws.Range("I99").Value = 397
$ws.Range("K99").Value = 397
$ws.Range("M99").Value = 1101
$ws.Range("H105").Value = 1419.5834
$ws.Range("I105").Value = 1419.5834
$ws.Range("K105").Value = 1419.5834
$ws.Range("M105").Value = 327.4166
$ws.Range("H107").Value = 5841.8335
$ws.Range("I107").Value = 5821
$ws.Range("K107").Value = 5821
$ws.Range("M107").Value = -3901

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 7097
$ws.Range("I62").Value = 7095
$ws.Range("J62").Value = 7099
$ws.Range("K62").Value = 7095
$ws.Range("L62").Value = 7099
$ws.Range("M62").Value = -6471
$ws.Range("N62").Value = -8347
$ws.Range("H65").Value = 7097
$ws.Range("I65").Value = 7095
$ws.Range("J65").Value = 7099
$ws.Range("K65").Value = 35475
$ws.Range("L65").Value = 35495
$ws.Range("M65").Value = -32355
$ws.Range("N65").Value = -41735
$ws.Range("H122").Value = 2599.5
$ws.Range("I122").Value = 3199
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 9597
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -7147
$ws.Range("N122").Value = -10900
$ws.Range("H134").Value = 5080.5
$ws.Range("I134").Value = 5080.5
$ws.Range("K134").Value = 15241.5
$ws.Range("M134").Value = -12706.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 909
$ws.Range("J22").Value = 1810
$ws.Range("L22").Value = 5430
$ws.Range("N22").Value = -5768
$ws.Range("H27").Value = 909
$ws.Range("J27").Value = 1810
$ws.Range("L27").Value = 5430
$ws.Range("N27").Value = -5634
$ws.Range("H113").Value = 800.8889
$ws.Range("I113").Value = 645.1667
$ws.Range("J113").Value = 1112.3334
$ws.Range("K113").Value = 1935.5001
$ws.Range("L113").Value = 3337.0002
$ws.Range("M113").Value = 234.4999
$ws.Range("N113").Value = -7677.0002
$ws.Range("H117").Value = 1167.5555
$ws.Range("I117").Value = 1495.8
$ws.Range("K117").Value = 4487.4
$ws.Range("M117").Value = -1045.4
$ws.Range("H121").Value = 591.1667
$ws.Range("I121").Value = 591.1667
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 1773.5001
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = -463.5001
$ws.Range("N121").ClearContents()
$ws.Range("H129").Value = 2288.5715
$ws.Range("I129").Value = 990
$ws.Range("J129").Value = 2505
$ws.Range("K129").Value = 2970
$ws.Range("L129").Value = 7515
$ws.Range("M129").Value = 2030
$ws.Range("N129").Value = -17515
$ws.Range("H131").Value = 1187.25
$ws.Range("I131").Value = 1187.25
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 3561.75
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 1478.25
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1929.5
$ws.Range("I97").Value = 315.5
$ws.Range("K97").Value = 315.5
$ws.Range("M97").Value = 180.5
$ws.Range("H107").Value = 438
$ws.Range("I107").Value = 439
$ws.Range("K107").Value = 439
$ws.Range("M107").Value = 1481
$ws.Range("H132").Value = 5332.3335
$ws.Range("I132").Value = 4666
$ws.Range("K132").Value = 13998
$ws.Range("M132").Value = -11468

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1829.6666
$ws.Range("J82").Value = 1490
$ws.Range("L82").Value = 1490
$ws.Range("N82").Value = -2212
$ws.Range("H85").Value = 1829.6666
$ws.Range("J85").Value = 1490
$ws.Range("L85").Value = 1490
$ws.Range("N85").Value = -3986
$ws.Range("H100").Value = 1992
$ws.Range("I100").Value = 1989
$ws.Range("K100").Value = 1989
$ws.Range("M100").Value = -1448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3180
$ws.Range("I96").Value = 3180
$ws.Range("K96").Value = 3180
$ws.Range("M96").Value = -1807
